$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: split/expand the "effort" paragraph to add the new EDSM/DJFMP
# sentence about tow_volume filtering, and tidy the original two-sentence
# run back into its own separate runs.
# ---------------------------------------------------------------------------
$find = $d.Content
$ok = $find.Find.Execute("To standardize catch across trawls, we calculated metrics of effort for each trawl. Effort was calculated following methods from the FMWT, Bay Study, and Suisun Study metadata and personal communications from the PIs. Sampling effort was quantified following the methods of the component surveys, as either a tow area (for the otter trawls) or a tow volume (for all other gear types). ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $ok) {
    throw "Could not locate target paragraph text for edit 1"
}

$find = $find.Paragraphs(1).Range

$paraXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
'<w:p w14:paraId="42570FEA" w14:textId="423BB0B0" w:rsidR="00410D7E" w:rsidRPr="00C80130" w:rsidRDefault="00410D7E" w:rsidP="005E02E2">' +
'<w:pPr><w:ind w:firstLine="720"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr>' +
'<w:r w:rsidRPr="00C80130"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>To standardize catch across trawls, we calculated metrics of effort for each trawl.</w:t></w:r>' +
'<w:r w:rsidRPr="00C80130"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
'<w:r w:rsidR="005E02E2" w:rsidRPr="00C80130"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">Sampling effort was quantified following the methods of the component surveys, as either a tow area (for the otter trawls) or a tow volume (for all other gear types). </w:t></w:r>' +
'<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">For EDSM and DJFMP data, </w:t></w:r>' +
'<w:proofErr w:type="spellStart"/>' +
'<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>tow_volume</w:t></w:r>' +
'<w:proofErr w:type="spellEnd"/>' +
'<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> was set to NA when debris was detected in the flowmeter, and samples were excluded with gear condition codes 3 (poor sampling), 4 (</w:t></w:r>' +
'<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>n</w:t></w:r>' +
'<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">o sample attempted), or 9 (fish captured outside of live box or cod end and could not be assigned to a specific tow). </w:t></w:r>' +
'</w:p>' +
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$find.InsertXML($paraXml)

# ---------------------------------------------------------------------------
# Edit 2: move <w:lastRenderedPageBreak/> from the "where FL = ..." run up
# to the start of the oMath run containing "FL=a+b×SL".
# ---------------------------------------------------------------------------

# The oMath equation paragraph isn't reachable by text search (Office Math
# runs don't expose readable plain text through Range.Text), so locate it
# via paragraph navigation from the preceding, plain-text paragraph.
$anchor = $d.Content
$ok2 = $anchor.Find.Execute("The equations are of the form", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok2) {
    throw "Could not locate paragraph preceding the oMath equation for edit 2"
}

$formPara = $anchor.Paragraphs(1)
$mathPara = $formPara.Next()
$wherePara = $mathPara.Next()

if ($wherePara.Range.Text -notlike "where FL*") {
    throw "Unexpected paragraph located for edit 2"
}

# 2a. Add <w:lastRenderedPageBreak/> to the oMath run.
$mathXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body>' +
'<w:p w14:paraId="41D00540" w14:textId="74023AD3" w:rsidR="00AE45C8" w:rsidRPr="00AE45C8" w:rsidRDefault="00AE45C8" w:rsidP="00AE45C8">' +
'<w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr>' +
'<m:oMathPara><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><w:lastRenderedPageBreak/><m:t>FL=a+b' + [char]0x00D7 + 'SL</m:t></m:r></m:oMath></m:oMathPara>' +
'</w:p>' +
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$mathPara.Range.InsertXML($mathXml)

# 2b. Remove <w:lastRenderedPageBreak/> from the "where FL = ..." run.
$whereXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
'<w:p w14:paraId="1497AF82" w14:textId="2B0A99DB" w:rsidR="00AE45C8" w:rsidRDefault="00AE45C8" w:rsidP="00AE45C8">' +
'<w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr>' +
'<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">where FL = fork length (or total length for species with no fork) in mm, a = intercept, b = slope, and SL = standard length in mm. </w:t></w:r>' +
'</w:p>' +
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$wherePara.Range.InsertXML($whereXml)

Write-Output "done"
